# Inventory change stock automation updated
# A1 keeps its displayed text ("NoOfSku") but B3's SKU label text changes
# from "test SKU 999" to "test SKU 992".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "test SKU 992"

# Move the active cell selection to B3
$ws.Range("B3").Select()
